# Update header row (row 1) labels so Power BI can automatically use the
# first row as a header when importing these tables.
#
# Sheets "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)" and "Emissoes Totais (MtCO2eq)" use the prefix
# "Ano " in front of each year.
#
# Sheet "Potencia Incremental - SIN(MW)" uses the prefix "Intervalo " in
# front of each year/interval value.
#
# Sheet "Custo Total (bilhões de R$)" only has a single year column (B1)
# and also gets the "Ano " prefix.

$wb = $excel.ActiveWorkbook

$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

$wsIncremental = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIncremental.Range("B1").Value = "Intervalo 2015"
$wsIncremental.Range("C1").Value = "Intervalo 2015-2030"
$wsIncremental.Range("D1").Value = "Intervalo 2031-2040"
$wsIncremental.Range("E1").Value = "Intervalo 2041-2050"

$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
